# Aula 12 - Algoritmos e Complexidade
#
# The running subtitle placeholder ("Algoritmos e Complexidade ...") that
# appears on the section-divider slides used to live in a single run for
# "Algoritmos e Complexidade" followed by a second run for " ...". The fix
# splits the first run into "Algoritmos e " / "Complexidade" and pads the
# trailing ellipsis run with a trailing space (" ... ") - exactly the
# content/run-shape change captured by the diff. Apply it wherever that
# exact text shows up (the text is duplicated, verbatim, across the
# section-divider slides), since the fix is purely mechanical/textual.

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text

        $oldFirst = "Algoritmos e Complexidade"
        $oldTail = " ..."

        if ($fullText -eq ($oldFirst + $oldTail)) {
            # Split "Algoritmos e Complexidade" into "Algoritmos e " + "Complexidade",
            # keeping each half's original run formatting (Times New Roman, sz 2000, ...).
            $part1 = $tr.Characters(1, 13)
            $part1.Text = "Algoritmos e "

            $part2 = $tr.Characters(14, 12)
            $part2.Text = "Complexidade"

            # Pad the trailing " ..." run with a trailing space -> " ... ".
            $part3 = $tr.Characters(26, 4)
            $part3.Text = " ... "
        }
    }
}
